# Natmi re-run (L1cam-Itga5, YoungD0) following Dr Hou's advice: the
# sending/target clusters now include "FAPs" alongside "ECs"/"sCs", and all
# expression-weight/specificity statistics were recomputed, growing the
# table from 6 data rows (A2:T7) to 9 data rows (A2:T10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "L1cam"
$ws.Cells.Item(2,3).Value = "Itga5"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 19.72083766666667
$ws.Cells.Item(2,8).Value = 59.162513
$ws.Cells.Item(2,9).Value = 0.8016210077351786
$ws.Cells.Item(2,10).Value = 0.8016210077351787
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 28.72417333333333
$ws.Cells.Item(2,14).Value = 86.17251999999999
$ws.Cells.Item(2,15).Value = 0.4233259107972328
$ws.Cells.Item(2,16).Value = 0.4233259107972328
$ws.Cells.Item(2,17).Value = 566.4647594158622
$ws.Cells.Item(2,18).Value = 5098.18283474276
$ws.Cells.Item(2,19).Value = 0.3393469432136901
$ws.Cells.Item(2,20).Value = 0.3393469432136901
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "L1cam"
$ws.Cells.Item(3,3).Value = "Itga5"
$ws.Cells.Item(3,4).Value = "sCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 19.72083766666667
$ws.Cells.Item(3,8).Value = 59.162513
$ws.Cells.Item(3,9).Value = 0.8016210077351786
$ws.Cells.Item(3,10).Value = 0.8016210077351787
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 30.56986233333333
$ws.Cells.Item(3,14).Value = 91.709587
$ws.Cells.Item(3,15).Value = 0.4505269713084062
$ws.Cells.Item(3,16).Value = 0.4505269713084062
$ws.Cells.Item(3,17).Value = 602.8632925680146
$ws.Cells.Item(3,18).Value = 5425.769633112131
$ws.Cells.Item(3,19).Value = 0.3611518847521225
$ws.Cells.Item(3,20).Value = 0.3611518847521226
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "L1cam"
$ws.Cells.Item(4,3).Value = "Itga5"
$ws.Cells.Item(4,4).Value = "FAPs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 19.72083766666667
$ws.Cells.Item(4,8).Value = 59.162513
$ws.Cells.Item(4,9).Value = 0.8016210077351786
$ws.Cells.Item(4,10).Value = 0.8016210077351787
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 8.559531999999999
$ws.Cells.Item(4,14).Value = 25.678596
$ws.Cells.Item(4,15).Value = 0.126147117894361
$ws.Cells.Item(4,16).Value = 0.126147117894361
$ws.Cells.Item(4,17).Value = 168.8011410746387
$ws.Cells.Item(4,18).Value = 1519.210269671748
$ws.Cells.Item(4,19).Value = 0.101122179769366
$ws.Cells.Item(4,20).Value = 0.101122179769366
$ws.Cells.Item(5,1).Value = "sCs"
$ws.Cells.Item(5,2).Value = "L1cam"
$ws.Cells.Item(5,3).Value = "Itga5"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.099159
$ws.Cells.Item(5,8).Value = 0.297477
$ws.Cells.Item(5,9).Value = 0.004030657259573097
$ws.Cells.Item(5,10).Value = 0.004030657259573097
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 28.72417333333333
$ws.Cells.Item(5,14).Value = 86.17251999999999
$ws.Cells.Item(5,15).Value = 0.4233259107972328
$ws.Cells.Item(5,16).Value = 0.4233259107972328
$ws.Cells.Item(5,17).Value = 2.84826030356
$ws.Cells.Item(5,18).Value = 25.63434273204
$ws.Cells.Item(5,19).Value = 0.001706281655520259
$ws.Cells.Item(5,20).Value = 0.00170628165552026
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "L1cam"
$ws.Cells.Item(6,3).Value = "Itga5"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.099159
$ws.Cells.Item(6,8).Value = 0.297477
$ws.Cells.Item(6,9).Value = 0.004030657259573097
$ws.Cells.Item(6,10).Value = 0.004030657259573097
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 30.56986233333333
$ws.Cells.Item(6,14).Value = 91.709587
$ws.Cells.Item(6,15).Value = 0.4505269713084062
$ws.Cells.Item(6,16).Value = 0.4505269713084062
$ws.Cells.Item(6,17).Value = 3.031276979111
$ws.Cells.Item(6,18).Value = 27.281492811999
$ws.Cells.Item(6,19).Value = 0.001815919807537708
$ws.Cells.Item(6,20).Value = 0.001815919807537708
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "L1cam"
$ws.Cells.Item(7,3).Value = "Itga5"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.099159
$ws.Cells.Item(7,8).Value = 0.297477
$ws.Cells.Item(7,9).Value = 0.004030657259573097
$ws.Cells.Item(7,10).Value = 0.004030657259573097
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 8.559531999999999
$ws.Cells.Item(7,14).Value = 25.678596
$ws.Cells.Item(7,15).Value = 0.126147117894361
$ws.Cells.Item(7,16).Value = 0.126147117894361
$ws.Cells.Item(7,17).Value = 0.8487546335879999
$ws.Cells.Item(7,18).Value = 7.638791702291999
$ws.Cells.Item(7,19).Value = 0.0005084557965151293
$ws.Cells.Item(7,20).Value = 0.0005084557965151295
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "L1cam"
$ws.Cells.Item(8,3).Value = "Itga5"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 4.781202
$ws.Cells.Item(8,8).Value = 14.343606
$ws.Cells.Item(8,9).Value = 0.1943483350052483
$ws.Cells.Item(8,10).Value = 0.1943483350052483
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 28.72417333333333
$ws.Cells.Item(8,14).Value = 86.17251999999999
$ws.Cells.Item(8,15).Value = 0.4233259107972328
$ws.Cells.Item(8,16).Value = 0.4233259107972328
$ws.Cells.Item(8,17).Value = 137.33607498968
$ws.Cells.Item(8,18).Value = 1236.02467490712
$ws.Cells.Item(8,19).Value = 0.08227268592802243
$ws.Cells.Item(8,20).Value = 0.08227268592802245
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "L1cam"
$ws.Cells.Item(9,3).Value = "Itga5"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 4.781202
$ws.Cells.Item(9,8).Value = 14.343606
$ws.Cells.Item(9,9).Value = 0.1943483350052483
$ws.Cells.Item(9,10).Value = 0.1943483350052483
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 30.56986233333333
$ws.Cells.Item(9,14).Value = 91.709587
$ws.Cells.Item(9,15).Value = 0.4505269713084062
$ws.Cells.Item(9,16).Value = 0.4505269713084062
$ws.Cells.Item(9,17).Value = 146.160686927858
$ws.Cells.Item(9,18).Value = 1315.446182350722
$ws.Cells.Item(9,19).Value = 0.087559166748746
$ws.Cells.Item(9,20).Value = 0.08755916674874602
$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,2).Value = "L1cam"
$ws.Cells.Item(10,3).Value = "Itga5"
$ws.Cells.Item(10,4).Value = "FAPs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 4.781202
$ws.Cells.Item(10,8).Value = 14.343606
$ws.Cells.Item(10,9).Value = 0.1943483350052483
$ws.Cells.Item(10,10).Value = 0.1943483350052483
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 8.559531999999999
$ws.Cells.Item(10,14).Value = 25.678596
$ws.Cells.Item(10,15).Value = 0.126147117894361
$ws.Cells.Item(10,16).Value = 0.126147117894361
$ws.Cells.Item(10,17).Value = 40.924851517464
$ws.Cells.Item(10,18).Value = 368.323663657176
$ws.Cells.Item(10,19).Value = 0.02451648232847982
$ws.Cells.Item(10,20).Value = 0.02451648232847982
